$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7532.6665
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7532.6665
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 22597.9995
$ws.Range("N69").Value = -24345.9995
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 7532.6665
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 7532.6665
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 67793.9985
$ws.Range("N72").Value = -76529.9985
$ws.Range("M72").ClearContents()

$ws.Range("H86").Value = 2963.158
$ws.Range("I86").Value = 2058.8333
$ws.Range("J86").Value = 4513.4287
$ws.Range("K86").Value = 2058.8333
$ws.Range("L86").Value = 4513.4287
$ws.Range("M86").Value = -935.8332999999998
$ws.Range("N86").Value = -6759.4287

$ws.Range("H89").Value = 2963.158
$ws.Range("I89").Value = 2058.8333
$ws.Range("J89").Value = 4513.4287
$ws.Range("K89").Value = 10294.1665
$ws.Range("L89").Value = 22567.1435
$ws.Range("M89").Value = -4678.166499999999
$ws.Range("N89").Value = -33799.14350000001

$ws.Range("H99").Value = 500.57144
$ws.Range("I99").Value = 500.57144
$ws.Range("K99").Value = 1501.71432
$ws.Range("M99").Value = -3.714320000000043

$ws.Range("H113").Value = 10325.3545
$ws.Range("I113").Value = 7854.1
$ws.Range("K113").Value = 7854.1
$ws.Range("M113").Value = -4600.1

$ws.Range("H118").Value = 3497.625
$ws.Range("I118").Value = 996.3333
$ws.Range("J118").Value = 4998.4
$ws.Range("K118").Value = 2988.9999
$ws.Range("L118").Value = 14995.2
$ws.Range("M118").Value = -1331.9999
$ws.Range("N118").Value = -18309.2

$ws.Range("H121").Value = 2026.2
$ws.Range("J121").Value = 2026.2
$ws.Range("L121").Value = 6078.6
$ws.Range("N121").Value = -9572.6

$ws.Range("H138").Value = 6253016.5
$ws.Range("I138").Value = 1592.6666
$ws.Range("J138").Value = 7356209
$ws.Range("K138").Value = 4777.9998
$ws.Range("L138").Value = 22068627
$ws.Range("M138").Value = 362.0002000000004
$ws.Range("N138").Value = -22078907

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7441.7285
$ws.Range("I32").Value = 3504.261
$ws.Range("K32").Value = 3504.261
$ws.Range("M32").Value = -3217.261

$ws.Range("H61").Value = 4857.2925
$ws.Range("I61").Value = 3841.611
$ws.Range("K61").Value = 3841.611
$ws.Range("M61").Value = -3629.611

$ws.Range("H74").Value = 50226.5
$ws.Range("I74").Value = 73213.78999999999
$ws.Range("K74").Value = 73213.78999999999
$ws.Range("M74").Value = -72339.78999999999

$ws.Range("H77").Value = 50226.5
$ws.Range("I77").Value = 73213.78999999999
$ws.Range("K77").Value = 366068.95
$ws.Range("M77").Value = -361700.95

$ws.Range("H132").Value = 4019.5854
$ws.Range("I132").Value = 3841.25
$ws.Range("K132").Value = 11523.75
$ws.Range("M132").Value = -8993.75

$ws.Range("H136").Value = 4857.2925
$ws.Range("I136").Value = 3841.611
$ws.Range("K136").Value = 11524.833
$ws.Range("M136").Value = -8974.832999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 25199.2
$ws.Range("J2").Value = 25199.2
$ws.Range("L2").Value = 25199.2
$ws.Range("N2").Value = -25425.2

$ws.Range("H105").Value = 2464.7908
$ws.Range("I105").Value = 2487.4285
$ws.Range("J105").Value = 2365.75
$ws.Range("K105").Value = 2487.4285
$ws.Range("L105").Value = 2365.75
$ws.Range("M105").Value = -740.4285
$ws.Range("N105").Value = -5859.75

$ws.Range("H134").Value = 4974.8335
$ws.Range("I134").Value = 3924.5
$ws.Range("J134").Value = 5500
$ws.Range("K134").Value = 11773.5
$ws.Range("L134").Value = 16500
$ws.Range("M134").Value = -9238.5
$ws.Range("N134").Value = -21570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 799.7143
$ws.Range("I22").Value = 291.5
$ws.Range("K22").Value = 291.5
$ws.Range("M22").Value = 58.5

$ws.Range("H99").Value = 3083.0833
$ws.Range("I99").Value = 2780.375
$ws.Range("J99").Value = 3688.5
$ws.Range("K99").Value = 2780.375
$ws.Range("L99").Value = 3688.5
$ws.Range("M99").Value = -1282.375
$ws.Range("N99").Value = -6684.5

$ws.Range("H105").Value = 1118.8
$ws.Range("I105").Value = 1367.4445
$ws.Range("K105").Value = 1367.4445
$ws.Range("M105").Value = 379.5554999999999

$ws.Range("H124").Value = 48100
$ws.Range("J124").Value = 48100
$ws.Range("L124").Value = 48100
$ws.Range("N124").Value = -53010

$ws.Range("H126").Value = 3083.0833
$ws.Range("I126").Value = 2780.375
$ws.Range("J126").Value = 3688.5
$ws.Range("K126").Value = 8341.125
$ws.Range("L126").Value = 11065.5
$ws.Range("M126").Value = -5871.125
$ws.Range("N126").Value = -16005.5

$ws.Range("H132").Value = 4261
$ws.Range("I132").Value = 3339.2666
$ws.Range("K132").Value = 10017.7998
$ws.Range("M132").Value = -7487.799800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 900
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H82").Value = 26796.125
$ws.Range("I82").Value = 15728.167
$ws.Range("K82").Value = 47184.501
$ws.Range("M82").Value = -46778.501

$ws.Range("H83").Value = 900
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H85").Value = 26796.125
$ws.Range("I85").Value = 15728.167
$ws.Range("K85").Value = 47184.501
$ws.Range("M85").Value = -45780.501

$ws.Range("H124").Value = 41211.727
$ws.Range("I124").Value = 29
$ws.Range("K124").Value = 87
$ws.Range("M124").Value = 4823

$ws.Range("H126").Value = 3143.3333
$ws.Range("I126").Value = 3030
$ws.Range("J126").Value = 3200
$ws.Range("K126").Value = 9090
$ws.Range("L126").Value = 9600
$ws.Range("M126").Value = -4150
$ws.Range("N126").Value = -19480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 3079
$ws.Range("I19").Value = 2800
$ws.Range("K19").Value = 2800
$ws.Range("M19").Value = -2512

$ws.Range("H102").Value = 58824430
$ws.Range("J102").Value = 250000800
$ws.Range("L102").Value = 250000800
$ws.Range("N102").Value = -250004044

$ws.Range("H126").Value = 21250.75
$ws.Range("I126").Value = 24631
$ws.Range("J126").Value = 4349.5
$ws.Range("K126").Value = 73893
$ws.Range("L126").Value = 13048.5
$ws.Range("M126").Value = -71423
$ws.Range("N126").Value = -17988.5

$ws.Range("H132").Value = 3510.195
$ws.Range("I132").Value = 2937.7585
$ws.Range("J132").Value = 4893.5835
$ws.Range("K132").Value = 8813.2755
$ws.Range("L132").Value = 14680.7505
$ws.Range("M132").Value = -6283.2755
$ws.Range("N132").Value = -19740.7505

$ws.Range("H136").Value = 23255.75
$ws.Range("J136").Value = 23255.75
$ws.Range("L136").Value = 69767.25
$ws.Range("N136").Value = -74867.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5921.643
$ws.Range("J40").Value = 6609
$ws.Range("L40").Value = 6609
$ws.Range("N40").Value = -6881

$ws.Range("H46").Value = 1104.375
$ws.Range("J46").Value = 1207.2
$ws.Range("L46").Value = 1207.2
$ws.Range("N46").Value = -1583.2

$ws.Range("H74").Value = 18750
$ws.Range("I74").Value = 17500
$ws.Range("K74").Value = 17500
$ws.Range("M74").Value = -16502

$ws.Range("H77").Value = 18750
$ws.Range("I77").Value = 17500
$ws.Range("K77").Value = 52500
$ws.Range("M77").Value = -47508

$ws.Range("H93").Value = 1988.7693
$ws.Range("I93").Value = 1995.909
$ws.Range("J93").Value = 1949.5
$ws.Range("K93").Value = 1995.909
$ws.Range("L93").Value = 1949.5
$ws.Range("M93").Value = -747.9090000000001
$ws.Range("N93").Value = -4445.5

$ws.Range("H136").Value = 3474.6875
$ws.Range("I136").Value = 3216.25
$ws.Range("K136").Value = 9648.75
$ws.Range("M136").Value = -7098.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2500.1667
$ws.Range("I136").Value = 1532.875
$ws.Range("K136").Value = 4598.625
$ws.Range("M136").Value = -2048.625
